$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Child")

$addresses = @{
    2  = "-5.2,-4.66"
    3  = "-0.07,9.35"
    4  = "7.33,-4.6"
    5  = "-0.99,0.61"
    6  = "6.22,-0.63"
    7  = "0.03,9.92"
    8  = "3.13,1.94"
    9  = "-4.29,-7.75"
    10 = "-1.97,-7.93"
    11 = "8.33,8.63"
    12 = "-1.87,4.17"
    13 = "-8.86,-5.32"
    14 = "8.79,0.28"
    15 = "2.1,3.74"
    16 = "5.04,3.77"
    17 = "9.82,7.5"
    18 = "-4.8,-6.74"
    19 = "6.07,9.47"
    20 = "9.89,-0.94"
    21 = "2.72,-7.47"
    22 = "0.86,9.91"
}

foreach ($row in $addresses.Keys) {
    $ws.Range("D$row").Value = $addresses[$row]
}
